# Add two new columns "I0" (I) and "IF" (J) to the right of the existing
# data, matching header style/format used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy formatting from the existing "IP" header (H1)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-55: literal values for the new I0 / IF columns.
$data = @(
    @(7, 8),
    @(5, 5),
    @(4, 4),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(10, 10),
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(10, 11),
    @(9, 9),
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
